$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: remove empty F10 cell, correct G10 district name
$ws.Range("F10").ClearContents()
$ws.Range("G10").Value = "Udupi"

# Row 16: correct G16 district name
$ws.Range("G16").Value = "Udupi"

# Rows 33-47 and 49-55: standardize district name to official name
$rows = @(33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,49,50,51,52,53,54,55)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = "Uttara Kannada (Karwar)"
}

# Row 48: remove empty F48 cell
$ws.Range("F48").ClearContents()
